# Apply cryptocurrency price/volume updates to match the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '23.768.06'
$ws.Cells.Item(2, 5).Value = '  +2.28%  '

$ws.Cells.Item(3, 4).Value = '1.656.42'
$ws.Cells.Item(3, 5).Value = '  +3.19%  '

$ws.Cells.Item(4, 4).Value = '''0.9980'
$ws.Cells.Item(4, 5).Value = '  -0.27%  '

$ws.Cells.Item(5, 5).Value = '  +0.76%  '

$ws.Cells.Item(6, 4).Value = '''0.9980'
$ws.Cells.Item(6, 5).Value = '  -0.28%  '

$ws.Cells.Item(7, 4).Value = '''0.3783'
$ws.Cells.Item(7, 5).Value = '  +0.56%  '

$ws.Cells.Item(8, 4).Value = '''52.89'
$ws.Cells.Item(8, 5).Value = '  +0.97%  '

$ws.Cells.Item(9, 4).Value = '''0.3687'
$ws.Cells.Item(9, 5).Value = '  +1.60%  '

$ws.Cells.Item(10, 4).Value = '''1.279'
$ws.Cells.Item(10, 5).Value = '  +0.43%  '

$ws.Cells.Item(11, 4).Value = '''0.08183'
$ws.Cells.Item(11, 5).Value = '  +0.35%  '

$ws.Cells.Item(12, 4).Value = '''0.9982'
$ws.Cells.Item(12, 5).Value = '  -0.25%  '

$ws.Cells.Item(13, 4).Value = '''23.29'
$ws.Cells.Item(13, 5).Value = '  +1.59%  '

$ws.Cells.Item(14, 4).Value = '''6.743'
$ws.Cells.Item(14, 5).Value = '  +2.17%  '

$ws.Cells.Item(15, 4).Value = '''0.00001277'
$ws.Cells.Item(15, 5).Value = '  +2.08%  '

$ws.Cells.Item(16, 4).Value = '''7.436'
$ws.Cells.Item(16, 5).Value = '  +0.66%  '

$ws.Cells.Item(17, 4).Value = '1.656.70'
$ws.Cells.Item(17, 5).Value = '  +3.11%  '

$ws.Cells.Item(18, 4).Value = '''95.52'
$ws.Cells.Item(18, 5).Value = '  +1.75%  '

$ws.Cells.Item(19, 5).Value = '  +0.00%  '

$ws.Cells.Item(20, 4).Value = '''18.54'
$ws.Cells.Item(20, 5).Value = '  +2.14%  '

$ws.Cells.Item(21, 4).Value = '''6.633'
$ws.Cells.Item(21, 5).Value = '  +1.42%  '

$ws.Cells.Item(22, 4).Value = '''0.9978'
$ws.Cells.Item(22, 5).Value = '  -0.52%  '

$ws.Cells.Item(23, 4).Value = '23.782.20'
$ws.Cells.Item(23, 5).Value = '  +2.32%  '

$ws.Cells.Item(24, 4).Value = '''13.04'
$ws.Cells.Item(24, 5).Value = '  +0.85%  '

$ws.Cells.Item(25, 4).Value = '''3.256'
$ws.Cells.Item(25, 5).Value = '  +4.93%  '

$ws.Cells.Item(26, 4).Value = '''2.429'
$ws.Cells.Item(26, 5).Value = '  -0.88%  '

$ws.Cells.Item(27, 5).Value = '  +1.57%  '

$ws.Cells.Item(28, 4).Value = '''152.22'
$ws.Cells.Item(28, 5).Value = '  +1.50%  '

$ws.Cells.Item(29, 4).Value = '''5.328'
$ws.Cells.Item(29, 5).Value = '  +0.94%  '

$ws.Cells.Item(30, 4).Value = '''137.64'
$ws.Cells.Item(30, 5).Value = '  +1.91%  '

$ws.Cells.Item(31, 4).Value = '''2.318'
$ws.Cells.Item(31, 5).Value = '  -3.17%  '

$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(32, 4).Value = '''7.069'
$ws.Cells.Item(32, 5).Value = '  +4.67%  '

$ws.Cells.Item(33, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(33, 4).Value = '1.844.45'
$ws.Cells.Item(33, 5).Value = '  +3.56%  '

$ws.Cells.Item(34, 4).Value = '''11.09'
$ws.Cells.Item(34, 5).Value = '  +6.71%  '

$ws.Cells.Item(35, 4).Value = '''0.9841'
$ws.Cells.Item(35, 5).Value = '  +2.06%  '

$ws.Cells.Item(36, 4).Value = '''0.02906'
$ws.Cells.Item(36, 5).Value = '  +5.30%  '

$ws.Cells.Item(37, 4).Value = '''6.425'
$ws.Cells.Item(37, 5).Value = '  +5.02%  '

$ws.Cells.Item(38, 4).Value = '''0.2598'
$ws.Cells.Item(38, 5).Value = '  +3.24%  '

$ws.Cells.Item(39, 4).Value = '''0.07376'
$ws.Cells.Item(39, 5).Value = '  -1.51%  '

$ws.Cells.Item(40, 4).Value = '''0.08928'
$ws.Cells.Item(40, 5).Value = '  +1.45%  '

$ws.Cells.Item(41, 4).Value = '''0.7252'
$ws.Cells.Item(41, 5).Value = '  +2.19%  '

$ws.Cells.Item(42, 4).Value = '''1.389'
$ws.Cells.Item(42, 5).Value = '  -2.36%  '

$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(43, 4).Value = '''16.69'
$ws.Cells.Item(43, 5).Value = '  +5.68%  '

$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(44, 4).Value = '''12.80'
$ws.Cells.Item(44, 5).Value = '  +2.55%  '

$ws.Cells.Item(45, 4).Value = '''0.6682'
$ws.Cells.Item(45, 5).Value = '  +2.25%  '

$ws.Cells.Item(46, 4).Value = '''2.404'
$ws.Cells.Item(46, 5).Value = '  +3.13%  '

$ws.Cells.Item(47, 4).Value = '''4.039'
$ws.Cells.Item(47, 5).Value = '  +0.75%  '

$ws.Cells.Item(48, 5).Value = '  -0.27%  '

$ws.Cells.Item(49, 4).Value = '''0.08113'
$ws.Cells.Item(49, 5).Value = '  +2.22%  '

$ws.Cells.Item(50, 4).Value = '''1.238'
$ws.Cells.Item(50, 5).Value = '  +2.42%  '

$ws.Cells.Item(51, 4).Value = '''129.48'
